# feature: new record, related name
# Add a new "is_new_record" header column (T1), matching the formatting
# of the preceding header cell (S1), update the row height to auto, and
# move the viewport / selection the same way the author's session did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy S1's formatting onto the new T1 cell, then set its text. Using
# copy/paste-special (formats only) reuses the existing style record
# instead of minting a new one, matching how Excel behaves when a header
# row's style is extended one column to the right.
$ws.Range("S1").Copy()
$ws.Range("T1").PasteSpecial(-4122)
$ws.Range("T1").Value = "is_new_record"
$excel.CutCopyMode = $false

# The header row had an explicit custom height (ht="13"); re-fit it so it
# reverts to the sheet's default auto height, matching the authored edit.
$ws.Rows.Item(1).AutoFit() | Out-Null

# Match the author's scrolled viewport / active selection at save time.
$win = $excel.ActiveWindow
$win.ScrollColumn = 17
$win.ScrollRow = 1
$ws.Range("U4").Select()
